# The document has two headers (header1.xml / header2.xml) each holding an
# inline BTEC logo picture, and two footers (footer1.xml / footer2.xml) each
# holding an inline Pearson logo picture. The picture's wp:docPr/@name (and
# pic:cNvPr/@name) attributes need their numeric suffix swapped:
#   BTEC logo pictures:    image1.jpg -> image2.jpg
#   Pearson logo pictures: image2.png -> image1.png
# (the id= attributes and the embedded media files themselves are unchanged)

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTEC_Logo-Orange pictures, image1.jpg -> image2.jpg ---
for ($h = 1; $h -le 2; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# --- Footers: Pearson logo pictures, image2.png -> image1.png ---
for ($f = 1; $f -le 2; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
